$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has data in rows 4, 5, 6, 8 (row 7 is an empty gap row).
# Target layout shifts everything up by 3 rows: 4->1, 5->2, 6->3, (7->4 stays empty), 8->5.
# Deleting three empty rows at the top accomplishes exactly that shift while
# preserving every cell's content, type and relative row gap.
$ws.Rows.Item(1).Delete()
$ws.Rows.Item(1).Delete()
$ws.Rows.Item(1).Delete()

# Row 2 (former row 5) holds the numeric values 1,2,3,4. Give the first two
# cells an integer ("0") number format, distinguishing them from the other
# numeric cells which keep the default/general format.
$ws.Range("A2:B2").NumberFormat = "0"

# Update the active selection to C2, matching the new selection in the file.
$ws.Range("C2").Select() | Out-Null
